# Edit SkillCardData.xlsx / Sheet1:
#  - Rows 3-5 ("等级3"/"等级4"/"等级5" cards): append the "若未找到，则重抽本牌"
#    fallback clause to the upgrade-swap effect text, and grow the row height
#    to fit the now-longer wrapped text.
#  - Update the sheet view's scroll position / active selection to match the
#    author's final cursor position (topLeftCell A7->A4, selection H9->F5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3: "等级3" alternate-upgrade-area card ---------------------------
$ws.Range("E3").Value = "进入备选升级牌区时：如果玩家未拥有《等级2》，则从升级牌堆中选1张《等级2》替换本牌，若未找到，则重抽本牌。<br>`n被动：可以使用至多包含3张牌的堆叠。"

# --- Row 4: "等级4" alternate-upgrade-area card ---------------------------
$ws.Range("E4").Value = "进入备选升级牌区时：如果玩家未拥有《等级3》，则从升级牌堆中选1张《等级3》替换本牌，若未找到，则重抽本牌。<br>被动：可以使用至多包含4张牌的堆叠。"

# --- Row 5: "等级5" alternate-upgrade-area card ---------------------------
$ws.Range("E5").Value = "进入备选升级牌区时：如果玩家未拥有《等级4》，则从升级牌堆中选1张《等级4》替换本牌，若未找到，则重抽本牌。<br>被动：可以使用至多包含5张牌的堆叠。"

# The longer text needs a taller wrapped row to display fully.
$ws.Rows(3).RowHeight = 270.75
$ws.Rows(4).RowHeight = 270.75
$ws.Rows(5).RowHeight = 270.75

# --- Restore the view state (scroll position + active cell) --------------
[void]$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
try { $win.TopLeftCell = $ws.Range("A4") } catch { }
[void]$ws.Range("F5").Select()
